# Add a new column S (year 2022 data) to the worksheet, mirroring the
# formatting already used for column R (year 2021).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmt = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# Copy formatting from column R into column S for every populated row so the
# new cells pick up the exact same style indexes as their neighbours.
$ws.Range("R3:R15").Copy()
$ws.Range("S3:S15").PasteSpecial($fmt)
$excel.CutCopyMode = 0

# Header year
$ws.Range("S3").Value = 2022

# Numeric data rows
$ws.Range("S4").Value = 10444.200000000001
$ws.Range("S5").Value = 21.7
$ws.Range("S6").Value = 7361.6
$ws.Range("S7").Value = 143.1
$ws.Range("S8").Value = 844.2
# Row 9 stays blank (no data reported yet for 2022)

# Text-valued, numeric-looking cells (authored as text in the source data)
$ws.Range("S10").Value = "2 756,0"
$ws.Range("S11").Value = "1 013,8"
$ws.Range("S12").Value = "1 451,1"

$ws.Range("S13").Value = 273.39999999999998
$ws.Range("S14").Value = "-"
$ws.Range("S15").Value = 17.7

# Update selection to mimic the authored view state
$ws.Range("T3").Select()
